$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 4137.254628375277
$ws.Range("C2").Value = 5786.418940849923
$ws.Range("D2").Value = 13395.78542994367
